# GroceryInventory.xlsx edit script
# Summary of changes (per commit message / diff):
#  - Added a row-1 freeze pane (view was scrolled/rearranged) and moved the
#    selection to H24.
#  - Bacon (row 11) and Grapes (row 24) "Limit" (column H) raised 6 -> 9.
#  - Banana (row 19) gained a "Has Specialty" flag with the "nmatx" (limit)
#    specialty type and its variables/limit filled in, mirroring the existing
#    Lime (row 17) specialty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Data edits
# ---------------------------------------------------------------------------

# Bacon (row 11): Limit 6 -> 9
$ws.Range("H11").Value = 9

# Banana (row 19): turn on "Has Specialty", type "nmatx", and fill in the
# specialty limit/variables (matches the pattern used by Lime in row 17).
$ws.Range("F19").Value = $true
$ws.Range("G19").Value = "nmatx"
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0.5

# Grapes (row 24): Limit 6 -> 9
$ws.Range("H24").Value = 9

# ---------------------------------------------------------------------------
# 2) View changes: freeze the header row and scroll/select like the target
# ---------------------------------------------------------------------------

$ws.Activate()

# Select row 2 so that freezing captures a 1-row split (header frozen).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Scroll the (now frozen) window down/right to roughly match the target
# view state, then land the real selection on H24 as in the target file.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H24").Select()
